$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "This activity would take approximately 15 seconds because it involves a brief interaction similar to `"Welcome and give application number,`" which took 28.2143 seconds. However, `"Welcome and get queue ticket`" is likely simpler and faster since it doesn't involve providing an application number, just issuing a ticket. Considering potential delays due to technical infrastructure limitations and the need for staff to be trained to efficiently handle the ticketing system, a slightly shorter time than the application number process is reasonable. Additionally, customer interaction time is minimal, as it involves just handing over a ticket."

$ws.Range("C3").Value = "This activity would take approximately 15 seconds because it involves a brief verbal interaction similar to `"Applicant hands out his/her documents`" (15.1429 seconds). The process requires the customer service officer to ask for consent, which is a straightforward question, and the applicant to respond. Given the technical infrastructure limitations and potential need for clarification or repetition due to language barriers or understanding, this time accounts for minor delays or retries. Staff training would ensure that the request is communicated clearly and efficiently, minimizing the time required."

$ws.Range("C4").Value = "Requesting a FIN ID from the application is similar to the activity where the CSO inputs the application number into the digital system, which takes 33.5714 seconds. Considering the technical infrastructure limitations in Ethiopia, there might be additional time required for system processing or connectivity issues. Staff training needs and potential customer interaction time could also add to the duration. Additionally, potential delays or retries due to system errors or incorrect data entry could extend the time. Therefore, a realistic estimate for this activity would be approximately 45 seconds."

$ws.Range("C5").Value = "Scanning a FIN QR code is a relatively straightforward task, similar to inputting an application number or taking a picture, which takes around 32-33 seconds. However, considering the technical infrastructure limitations in Ethiopia, potential delays due to network issues, and the need for staff to be adequately trained to handle QR code scanning efficiently, an additional buffer time is necessary. Customer interaction time for explaining the process and potential retries if the scan fails initially also contribute to the estimated time. Therefore, a realistic processing time for scanning a FIN QR code would be approximately 45 seconds."
$ws.Range("D5").Value = 45

$ws.Range("C6").Value = "This activity would take approximately 65 seconds because it involves orally gathering information from the applicant and entering it into the DARIS system. Similar activities, such as filling in personal data and other information into the digital system, take 48.7857 and 93.0 seconds respectively. Considering the technical infrastructure limitations in Ethiopia, potential delays, and the need for staff to be adequately trained to handle oral data entry efficiently, a moderate time estimate is reasonable. Additionally, customer interaction time and the possibility of needing to repeat or clarify information contribute to this duration."
$ws.Range("D6").Value = 65

$ws.Range("C7").Value = "This activity would take approximately 65 seconds because printing tasks such as `"CSO prints digital application (3 copies)`" take 77.7857 seconds, indicating that printing infrastructure might be slower due to technical limitations. The `"Print PA certificate`" might involve fewer copies or simpler documents, reducing time slightly. Staff training and customer interaction are minimal for printing, but potential delays or retries due to printer malfunctions or paper jams could occur, justifying a time estimate slightly lower than the digital application print time."

$ws.Range("C8").Value = "Stamping a PA certificate is similar to the `"Stamp documents`" activity, which takes 14.4286 seconds. Given the technical infrastructure and staff training in Ethiopia, a slight increase in time might be necessary to account for potential delays or retries, especially if the certificate requires additional verification or handling. Therefore, estimating 15 seconds is realistic, considering the need for precision and possible customer interaction during the stamping process."

$ws.Range("C9").Value = "Signing a PA certificate is similar to the activity where the applicant signs on the 3 copies to confirm all data are correct, which takes 52.0714 seconds. Considering that the PA certificate might be a single document, the signing process could be quicker. However, taking into account potential delays due to customer interaction, staff guidance, and the need for accuracy, a realistic estimate would be slightly less than the time taken for signing multiple copies. Therefore, 45 seconds is a reasonable estimate, accounting for the technical infrastructure and staff training in Ethiopia."
